$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet
$ws.Name = "Through 2022-04-13"

# Update header cell I1 (shared string value)
$ws.Range("I1").Value = "2022 (through 04-13)"

# Update I5 (May, 2022 column): 48 -> 50
$ws.Range("I5").Value = 50

# Update I14 (Total, 2022 column): 482 -> 484
$ws.Range("I14").Value = 484
